# Generate Report for Handback
# Refresh the handback-status report timestamps for the file that was just
# handed back (298558b4-a09c-4d42-8ff8-7e1ef7710dc7), across the Overview
# sheet and each per-language detail sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" column (G) for the 298558b4 row (row 2)
$overview.Range("G2").Value = "2016-08-24 14:55:26"

# zh-cn detail sheet, row 2 (298558b4):
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$zhcn.Range("H2").Value = "2016-08-24 14:55:21"
$zhcn.Range("K2").Value = "2016-08-24 14:55:39"

# de-de detail sheet, row 2 (298558b4):
#   H = Correspond Handoff Datetime, K = Correspond Handback DateTime
$dede.Range("H2").Value = "2016-08-24 14:55:26"
$dede.Range("K2").Value = "2016-08-24 14:55:46"
